$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 549. This shifts the existing rows
# 549-638 down to 550-639, preserving all of their data/styles.
$ws.Rows.Item(549).Insert()

# Populate the newly inserted row 549 with the new record's data.
# Columns A,B,C,E,F,G,H,I,R are identical to the surrounding rows, so
# Excel's insert (which copies formatting/values from the row above)
# already leaves correct-looking values there in most implementations,
# but we set everything explicitly to be safe/explicit.
$ws.Cells.Item(549, 1).Value = 3
$ws.Cells.Item(549, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(549, 3).Value = "Coquimbo"
$ws.Cells.Item(549, 4).Value = 45180
$ws.Cells.Item(549, 5).Value = 5
$ws.Cells.Item(549, 6).Value = 100112043
$ws.Cells.Item(549, 7).Value = "Pepino ensalada"
$ws.Cells.Item(549, 8).Value = "Sin especificar"
$ws.Cells.Item(549, 9).Value = "Primera"
$ws.Cells.Item(549, 10).Value = 100
$ws.Cells.Item(549, 11).Value = 10000
$ws.Cells.Item(549, 12).Value = 11000
$ws.Cells.Item(549, 13).Value = 10500
$ws.Cells.Item(549, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(549, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(549, 16).Value = 175
$ws.Cells.Item(549, 17).Value = 60
$ws.Cells.Item(549, 18).Value = "Hortaliza"
